# Updates the cryptos list values (Price / Volume(1h)) per the latest scrape,
# mirroring the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13/14 swap place (Solana now ranks above WrappedEther) and rows 41/42
# swap place (TheSandbox now ranks above InternetComputer(DFINITY)), so we
# set B..E fully for those rows. For all the rest, only D (Price) and
# E (Volume 1h) change.

$rows = @(
    @{ Row = 2;  D = "28.059.33";   E = "  -0.15%  " },
    @{ Row = 3;  D = "1.906.28";    E = "  +2.16%  " },
    @{ Row = 4;  D = "1.002";       E = "  -0.14%  " },
    @{ Row = 5;  D = "312.76";      E = "  +0.18%  " },
    @{ Row = 6;  D = "1.002";       E = "  -0.17%  " },
    @{ Row = 7;  D = "0.5050";      E = "  +1.08%  " },
    @{ Row = 8;  D = "0.3920";      E = "  +0.21%  " },
    @{ Row = 9;  D = "0.09585";     E = "  -0.51%  " },
    @{ Row = 10; D = "1.134";       E = "  +0.33%  " },
    @{ Row = 11; D = "42.03";       E = "  +2.93%  " },
    @{ Row = 12; D = "6.379";       E = "  -1.02%  " },
    @{ Row = 13; B = "Solana";        C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol";                               D = "20.78";      E = "  -0.24%  " },
    @{ Row = 14; B = "WrappedEther";  C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth";                    D = "1.891.23";   E = "  +1.15%  " },
    @{ Row = 15; D = "1.002";       E = "  -0.17%  " },
    @{ Row = 16; D = "7.295";       E = "  -0.94%  " },
    @{ Row = 17; D = "0.00001115"; E = "  -0.93%  " },
    @{ Row = 18; D = "92.11";       E = "  -0.81%  " },
    @{ Row = 19; D = "0.06599";     E = "  -0.30%  " },
    @{ Row = 20; D = "17.83";       E = "  +2.63%  " },
    @{ Row = 21; D = "1.002";       E = "  -0.21%  " },
    @{ Row = 22; D = "6.210";       E = "  +1.33%  " },
    @{ Row = 23; D = "28.118.34";   E = "  -0.19%  " },
    @{ Row = 24; D = "11.26";       E = "  +0.04%  " },
    @{ Row = 25; D = "2.302";       E = "  +0.54%  " },
    @{ Row = 26; D = "2.665";       E = "  +4.74%  " },
    @{ Row = 27; D = "2.123.82";    E = "  +1.85%  " },
    @{ Row = 28; D = "20.79";       E = "  -1.41%  " },
    @{ Row = 29; E = "  -0.28%  " },
    @{ Row = 30; D = "127.02";      E = "  -0.32%  " },
    @{ Row = 31; D = "1.084";       E = "  +2.99%  " },
    @{ Row = 32; D = "0.1061";      E = "  +0.28%  " },
    @{ Row = 33; D = "5.612";       E = "  +0.21%  " },
    @{ Row = 34; D = "3.616";       E = "  +0.02%  " },
    @{ Row = 35; D = "9.624";       E = "  +2.02%  " },
    @{ Row = 36; D = "0.06602";     E = "  -2.03%  " },
    @{ Row = 37; D = "0.02424";     E = "  +1.65%  " },
    @{ Row = 38; D = "1.231";       E = "  +0.72%  " },
    @{ Row = 39; D = "0.2180";      E = "  +0.21%  " },
    @{ Row = 40; D = "1.276";       E = "  +8.86%  " },
    @{ Row = 41; B = "TheSandbox";               C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand";                               D = "0.6335"; E = "  +1.21%  " },
    @{ Row = 42; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";                    D = "4.981";  E = "  -0.17%  " },
    @{ Row = 43; D = "11.35";       E = "  -0.68%  " },
    @{ Row = 44; D = "1.002";       E = "  -0.06%  " },
    @{ Row = 45; D = "13.24";       E = "  -1.54%  " },
    @{ Row = 46; D = "0.5984";      E = "  +0.19%  " },
    @{ Row = 47; D = "3.725";       E = "  +1.50%  " },
    @{ Row = 48; D = "1.278";       E = "  +0.40%  " },
    @{ Row = 49; E = "  +2.07%  " },
    @{ Row = 50; D = "122.75";      E = "  -1.31%  " },
    @{ Row = 51; E = "  -0.86%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $item.C }
    # Prices are stored as plain text in the source sheet (e.g. "28.059.33",
    # "0.5050"), not real numbers, so force text with a leading apostrophe -
    # otherwise numeric-looking strings like "1.002" or "6.210" would be
    # coerced to the Number type 1.002 / 6.21, silently dropping their
    # trailing zeros / thousands-dot formatting.
    if ($item.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = "'" + $item.D }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $item.E }
}
